# courseleafPatchControl.xlsx: update git directives to include the target directory
$wb = $excel.ActiveWorkbook

# --- 1. Update the "target" column for the git-sourced rows so the copy ---
# --- destination includes the program directory / component folder.    ---

# "cat" sheet: courseleaf git row (row 3) -> /web/<progDir>
$wsCat = $wb.Worksheets.Item("cat")
$wsCat.Range("C3").Value = "/web/<progDir>"

# "cim" sheet: cim git row (row 3) -> /web/<progDir>/ciim
$wsCim = $wb.Worksheets.Item("cim")
$wsCim.Range("C3").Value = "/web/<progDir>/ciim"

# "cat" sheet: navmaster git row (row 6) -> /web/<progDir>/pdf
$wsCat.Range("C6").Value = "/web/<progDir>/pdf"

# "pdfgen" sheet: pdfgen git row (row 3) -> /web/<progDir>pdf
$wsPdfgen = $wb.Worksheets.Item("pdfgen")
$wsPdfgen.Range("C3").Value = "/web/<progDir>pdf"

# "formbuilder" sheet: formbuilder git row (row 3) -> /web/<progDir>/formbuilder
$wsFormbuilder = $wb.Worksheets.Item("formbuilder")
$wsFormbuilder.Range("C3").Value = "/web/<progDir>/formbuilder"

# --- 2. Move the selection cursor on the edited sheets to C3 (where the ---
# --- edit was made), and switch the active tab to "formbuilder".       ---
[void]$wsCat.Range("C3").Select()
[void]$wsCim.Range("C3").Select()
[void]$wsPdfgen.Range("C3").Select()

# Activating "formbuilder" last makes it the workbook's active tab, and
# selecting C3 there sets its cursor/selection to match.
[void]$wsFormbuilder.Activate()
[void]$wsFormbuilder.Range("C3").Select()
